$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("G2").Value = "2016-08-17 22:47:00"
$ws1.Range("G3").Value = "2016-08-17 22:47:00"

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("H3").Value = "2016-08-17 22:46:54"
$ws2.Range("K3").Value = "2016-08-17 22:47:25"

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("H3").Value = "2016-08-17 22:47:00"
$ws3.Range("K3").Value = "2016-08-17 22:47:33"
